$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 121, shifting existing rows 121:195 down to 122:196
$ws.Rows.Item(121).Insert()

# Copy the static (series-identifying) columns from the row now at 122 (old row 121)
$ws.Range("A121").Value = $ws.Range("A122").Value2
$ws.Range("B121").Value = $ws.Range("B122").Value2
$ws.Range("C121").Value = $ws.Range("C122").Value2
$ws.Range("E121").Value = $ws.Range("E122").Value2
$ws.Range("F121").Value = $ws.Range("F122").Value2
$ws.Range("G121").Value = $ws.Range("G122").Value2
$ws.Range("H121").Value = $ws.Range("H122").Value2
$ws.Range("I121").Value = $ws.Range("I122").Value2
$ws.Range("R121").Value = $ws.Range("R122").Value2

# New row's own data
$ws.Range("D121").Value = 44767
$ws.Range("D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J121").Value = 120
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 1200
$ws.Range("M121").Value = 1200
$ws.Range("N121").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O121").Value = "Región de Los Lagos"
$ws.Range("P121").Value = 800
$ws.Range("Q121").Value = 1.5
